$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update title text (October -> November) -------------------------
$ws.Range("A1").Value = "Table 1.2.D. Net Generation by Energy Source:  Industrial Sector, 2006-November 2016"

# --- 2. Insert a new row for "November" before the old "Year to Date" ---
#        block (old row 53), pushing everything below down by one.
$ws.Rows.Item(53).Insert()

# Copy formatting from the October data row (row 52, now still row 52)
# down into the freshly inserted row 53 so the styles match (s=9/s=10).
$ws.Range("A52:P52").Copy()
$ws.Range("A53:P53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Populate the new November row (row 53) ---------------------------
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 595
$ws.Range("C53").Value = 40
$ws.Range("D53").Value = 69
$ws.Range("E53").Value = 7782
$ws.Range("F53").Value = 641
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 68
$ws.Range("I53").Value = "NM"
$ws.Range("J53").Value = 2379
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 407
$ws.Range("M53").Value = 11983
$ws.Range("N53").Value = 123
$ws.Range("O53").Value = "NM"
$ws.Range("P53").Value = "NM"

# --- 4. Update the "Year to Date" block (rows shifted down to 55-57) ----
# Row 55: 2014
$ws.Range("A55").Value = 2014
$ws.Range("B55").Value = 11326
$ws.Range("C55").Value = 503
$ws.Range("D55").Value = 1268
$ws.Range("E55").Value = 78540
$ws.Range("F55").Value = 7874
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 1157
$ws.Range("I55").Value = 16
$ws.Range("J55").Value = 26149
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 4494
$ws.Range("M55").Value = 131326
$ws.Range("N55").Value = 1065
$ws.Range("O55").Value = 1081
$ws.Range("P55").Value = 1081

# Row 56: 2015
$ws.Range("A56").Value = 2015
$ws.Range("B56").Value = 10064
$ws.Range("C56").Value = 525
$ws.Range("D56").Value = 912
$ws.Range("E56").Value = 80218
$ws.Range("F56").Value = 8595
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 1265
$ws.Range("I56").Value = 20
$ws.Range("J56").Value = 26155
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 4988
$ws.Range("M56").Value = 132742
$ws.Range("N56").Value = 1358
$ws.Range("O56").Value = 1378
$ws.Range("P56").Value = 1378

# Row 57: 2016
$ws.Range("A57").Value = 2016
$ws.Range("B57").Value = 8544
$ws.Range("C57").Value = 459
$ws.Range("D57").Value = 871
$ws.Range("E57").Value = 84298
$ws.Range("F57").Value = 8248
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 1177
$ws.Range("I57").Value = 28
$ws.Range("J57").Value = 25751
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 4776
$ws.Range("M57").Value = 134153
$ws.Range("N57").Value = 1720
$ws.Range("O57").Value = 1748
$ws.Range("P57").Value = 1748

# --- 5. Update "Rolling 12 Months Ending in October" -> "...November" ---
#        (this merged header row is now row 58)
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# --- 6. Update the rolling-12-months data rows (now rows 59-60) ---------
# Row 59: 2015
$ws.Range("A59").Value = 2015
$ws.Range("B59").Value = 11079
$ws.Range("C59").Value = 567
$ws.Range("D59").Value = 1033
$ws.Range("E59").Value = 87887
$ws.Range("F59").Value = 9385
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 1389
$ws.Range("I59").Value = 21
$ws.Range("J59").Value = 28665
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 5472
$ws.Range("M59").Value = 145498
$ws.Range("N59").Value = 1432
$ws.Range("O59").Value = 1453
$ws.Range("P59").Value = 1453

# Row 60: 2016
$ws.Range("A60").Value = 2016
$ws.Range("B60").Value = 9377
$ws.Range("C60").Value = 497
$ws.Range("D60").Value = "NM"
$ws.Range("E60").Value = 92436
$ws.Range("F60").Value = 9054
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 1322
$ws.Range("I60").Value = "NM"
$ws.Range("J60").Value = 28210
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 5251
$ws.Range("M60").Value = 147124
$ws.Range("N60").Value = 1813
$ws.Range("O60").Value = "NM"
$ws.Range("P60").Value = "NM"
